# Auto-generated edit script: apply updated market-price data to leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3295.1428
$ws.Range("I106").Value = 3510
$ws.Range("K106").Value = 3510
$ws.Range("M106").Value = -2879
$ws.Range("H107").Value = 1700
$ws.Range("I107").Value = 1400
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1400
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 520
$ws.Range("N107").Value = -5840
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
$ws.Range("H110").Value = 33700
$ws.Range("J110").Value = 33700
$ws.Range("L110").Value = 33700
$ws.Range("N110").Value = -41880
$ws.Range("H116").Value = 1833.4615
$ws.Range("I116").Value = 1766.8182
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 1766.8182
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = 1675.1818
$ws.Range("N116").Value = -9084
$ws.Range("H123").Value = 71457.78
$ws.Range("J123").Value = 71457.78
$ws.Range("L123").Value = 71457.78
$ws.Range("N123").Value = -81257.78
$ws.Range("H138").Value = 5716949.5
$ws.Range("I138").Value = 2087.8823
$ws.Range("J138").Value = 11114319
$ws.Range("K138").Value = 6263.646900000001
$ws.Range("L138").Value = 33342957
$ws.Range("M138").Value = -1123.646900000001
$ws.Range("N138").Value = -33353237

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4807.1
$ws.Range("I32").Value = 4524.3896
$ws.Range("J32").Value = 10178.6
$ws.Range("K32").Value = 4524.3896
$ws.Range("L32").Value = 10178.6
$ws.Range("M32").Value = -4237.3896
$ws.Range("N32").Value = -10752.6
$ws.Range("H82").Value = 33000
$ws.Range("J82").Value = 33000
$ws.Range("L82").Value = 33000
$ws.Range("N82").Value = -33722
$ws.Range("H85").Value = 33000
$ws.Range("J85").Value = 33000
$ws.Range("L85").Value = 33000
$ws.Range("N85").Value = -35496
$ws.Range("H122").Value = 5052725.5
$ws.Range("I122").Value = 2039.7273
$ws.Range("J122").Value = 10103411
$ws.Range("K122").Value = 6119.1819
$ws.Range("L122").Value = 30310233
$ws.Range("M122").Value = -3669.1819
$ws.Range("N122").Value = -30315133
$ws.Range("H131").Value = 57583.332
$ws.Range("J131").Value = 57583.332
$ws.Range("L131").Value = 57583.332
$ws.Range("N131").Value = -67663.33199999999
$ws.Range("H138").Value = 74514.336
$ws.Range("J138").Value = 74514.336
$ws.Range("L138").Value = 74514.336
$ws.Range("N138").Value = -84794.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1972.6666
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 1972.6666
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 1972.6666
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -2466.6666
$ws.Range("H59").Value = 100000
$ws.Range("J59").Value = 100000
$ws.Range("L59").Value = 100000
$ws.Range("N59").Value = -101694

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 277724.34
$ws.Range("I31").Value = 74439.78999999999
$ws.Range("J31").Value = 372590.47
$ws.Range("K31").Value = 74439.78999999999
$ws.Range("L31").Value = 372590.47
$ws.Range("M31").Value = -74144.78999999999
$ws.Range("N31").Value = -373180.47
$ws.Range("H34").Value = 277724.34
$ws.Range("I34").Value = 74439.78999999999
$ws.Range("J34").Value = 372590.47
$ws.Range("K34").Value = 74439.78999999999
$ws.Range("L34").Value = 372590.47
$ws.Range("M34").Value = -74237.78999999999
$ws.Range("N34").Value = -372994.47
$ws.Range("H58").Value = 18869464
$ws.Range("I58").Value = 27028386
$ws.Range("J58").Value = 1958.625
$ws.Range("K58").Value = 27028386
$ws.Range("L58").Value = 1958.625
$ws.Range("M58").Value = -27028183
$ws.Range("N58").Value = -2364.625
$ws.Range("H136").Value = 18869464
$ws.Range("I136").Value = 27028386
$ws.Range("J136").Value = 1958.625
$ws.Range("K136").Value = 81085158
$ws.Range("L136").Value = 5875.875
$ws.Range("M136").Value = -81082608
$ws.Range("N136").Value = -10975.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 542.4231
$ws.Range("I5").Value = 445.75
$ws.Range("J5").Value = 1702.5
$ws.Range("K5").Value = 1337.25
$ws.Range("L5").Value = 5107.5
$ws.Range("M5").Value = -1225.25
$ws.Range("N5").Value = -5331.5
$ws.Range("H40").Value = 231.3
$ws.Range("I40").Value = 44.142857
$ws.Range("J40").Value = 668
$ws.Range("K40").Value = 176.571428
$ws.Range("L40").Value = 2672
$ws.Range("M40").Value = -107.571428
$ws.Range("N40").Value = -2810
$ws.Range("H80").Value = 1398.4286
$ws.Range("J80").Value = 1716.8334
$ws.Range("L80").Value = 5150.5002
$ws.Range("N80").Value = -7022.5002
$ws.Range("H83").Value = 1398.4286
$ws.Range("J83").Value = 1716.8334
$ws.Range("L83").Value = 15451.5006
$ws.Range("N83").Value = -24811.5006
$ws.Range("H114").Value = 9804503
$ws.Range("I114").Value = 475.58334
$ws.Range("J114").Value = 15152154
$ws.Range("K114").Value = 1426.75002
$ws.Range("L114").Value = 45456462
$ws.Range("M114").Value = 1827.24998
$ws.Range("N114").Value = -45462970
$ws.Range("H122").Value = 915.8333
$ws.Range("I122").Value = 213.66667
$ws.Range("J122").Value = 1266.9166
$ws.Range("K122").Value = 1923.00003
$ws.Range("L122").Value = 11402.2494
$ws.Range("M122").Value = 526.9999699999998
$ws.Range("N122").Value = -16302.2494
$ws.Range("H131").Value = 877.3461
$ws.Range("J131").Value = 1000.75
$ws.Range("L131").Value = 3002.25
$ws.Range("N131").Value = -13082.25
$ws.Range("H132").Value = 3662.5908
$ws.Range("J132").Value = 4424.3335
$ws.Range("L132").Value = 39819.0015
$ws.Range("N132").Value = -44879.0015
$ws.Range("H134").Value = 6181.6177
$ws.Range("I134").Value = 3398.077
$ws.Range("J134").Value = 7904.7617
$ws.Range("K134").Value = 10194.231
$ws.Range("L134").Value = 23714.2851
$ws.Range("M134").Value = -5124.231
$ws.Range("N134").Value = -33854.2851
$ws.Range("H135").Value = 542.4231
$ws.Range("I135").Value = 445.75
$ws.Range("J135").Value = 1702.5
$ws.Range("K135").Value = 4011.75
$ws.Range("L135").Value = 15322.5
$ws.Range("M135").Value = -1476.75
$ws.Range("N135").Value = -20392.5
$ws.Range("H137").Value = 1445.3636
$ws.Range("I137").Value = 901
$ws.Range("J137").Value = 3895
$ws.Range("K137").Value = 2703
$ws.Range("L137").Value = 11685
$ws.Range("M137").Value = 2397
$ws.Range("N137").Value = -21885
$ws.Range("H138").Value = 4259.091
$ws.Range("I138").Value = 1081.3334
$ws.Range("K138").Value = 3244.0002
$ws.Range("M138").Value = 1895.9998
$ws.Range("H139").Value = 4381.875
$ws.Range("I139").Value = 1874.2778
$ws.Range("J139").Value = 7605.9287
$ws.Range("K139").Value = 5622.8334
$ws.Range("L139").Value = 22817.7861
$ws.Range("M139").Value = -482.8334000000004
$ws.Range("N139").Value = -33097.7861
$ws.Range("H140").Value = 3053.9546
$ws.Range("I140").Value = 4682.353
$ws.Range("J140").Value = 2028.6666
$ws.Range("K140").Value = 14047.059
$ws.Range("L140").Value = 6085.9998
$ws.Range("M140").Value = -8867.059000000001
$ws.Range("N140").Value = -16445.9998
$ws.Range("H141").Value = 14338.158
$ws.Range("I141").Value = 5146.4287
$ws.Range("J141").Value = 19700
$ws.Range("K141").Value = 15439.2861
$ws.Range("L141").Value = 59100
$ws.Range("M141").Value = -10259.2861
$ws.Range("N141").Value = -69460

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 22942.666
$ws.Range("I44").Value = 8828
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 8828
$ws.Range("L44").Value = 30000
$ws.Range("M44").Value = -8232
$ws.Range("N44").Value = -31192
$ws.Range("H70").Value = 45200
$ws.Range("I70").Value = 67543.75
$ws.Range("J70").Value = 5477.778
$ws.Range("K70").Value = 67543.75
$ws.Range("L70").Value = 5477.778
$ws.Range("M70").Value = -67273.75
$ws.Range("N70").Value = -6017.778
$ws.Range("H73").Value = 45200
$ws.Range("I73").Value = 67543.75
$ws.Range("J73").Value = 5477.778
$ws.Range("K73").Value = 67543.75
$ws.Range("L73").Value = 5477.778
$ws.Range("M73").Value = -66607.75
$ws.Range("N73").Value = -7349.778
$ws.Range("H122").Value = 1651.5416
$ws.Range("I122").Value = 1242.8823
$ws.Range("J122").Value = 2644
$ws.Range("K122").Value = 3728.6469
$ws.Range("L122").Value = 7932
$ws.Range("M122").Value = -1278.6469
$ws.Range("N122").Value = -12832

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1163.6818
$ws.Range("I16").Value = 861.6111
$ws.Range("K16").Value = 861.6111
$ws.Range("M16").Value = -691.6111
$ws.Range("H87").Value = 44589
$ws.Range("J87").Value = 44589
$ws.Range("L87").Value = 44589
$ws.Range("N87").Value = -46835
$ws.Range("H90").Value = 44589
$ws.Range("J90").Value = 44589
$ws.Range("L90").Value = 133767
$ws.Range("N90").Value = -144999
